# Applies the "Comments, RD2, test cases fixed -Cooper" edit.
# Updates the 4 test-case rows (9-12) on Sheet1: swaps the test description
# and package-color columns, drops the now-unused "months" input column (C),
# and fixes the row-12 formula which mistakenly referenced B5 instead of B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Green package, 2 GB
$ws.Range("A9").Value = "Green Package 2 GB"
$ws.Range("B9").Value = "green"
$ws.Range("C9").ClearContents()

# Row 10 - Blue package, 2 GB
$ws.Range("A10").Value = "Blue Package 2 GB"
$ws.Range("B10").Value = "blue"
$ws.Range("C10").ClearContents()

# Row 11 - Purple package
$ws.Range("A11").Value = "Purple Package"
$ws.Range("B11").Value = "purple"
$ws.Range("C11").ClearContents()

# Row 12 - Green package, 5 GB (also fix the stray B5 reference -> B2)
$ws.Range("A12").Value = "Green Package 5 GB"
$ws.Range("B12").Value = "green"
$ws.Range("C12").ClearContents()
$ws.Range("F12").Formula = "=IF(D12>B2,49.99+15*(D12-B2),49.99)"

# Header row - the "input: months" column is no longer used
$ws.Range("C8").ClearContents()

# Restore the active selection to F13
$ws.Range("F13").Select()
